# Weekly price-sheet update: a new (most recent) price observation is
# inserted as the new row 507, pushing the existing rows 507-589 down to
# 508-590. The sheet's used range grows from A1:R589 to A1:R590.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 507 - this shifts rows 507:589 down
# to 508:590 and extends the sheet dimension accordingly.
$ws.Rows.Item(507).Insert()

# Populate the newly inserted row 507 with the new weekly data point.
$ws.Cells.Item(507, 1).Value = 4
$ws.Cells.Item(507, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(507, 3).Value = "Los Lagos"
$ws.Cells.Item(507, 4).Value = 44951
$ws.Cells.Item(507, 5).Value = 10
$ws.Cells.Item(507, 6).Value = 100114001
$ws.Cells.Item(507, 7).Value = "Papa"
$ws.Cells.Item(507, 8).Value = "Patagonia"
$ws.Cells.Item(507, 9).Value = "1a nueva(o)"
$ws.Cells.Item(507, 10).Value = 300
$ws.Cells.Item(507, 11).Value = 13000
$ws.Cells.Item(507, 12).Value = 13000
$ws.Cells.Item(507, 13).Value = 13000
$ws.Cells.Item(507, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(507, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(507, 16).Value = 520
$ws.Cells.Item(507, 17).Value = 25
$ws.Cells.Item(507, 18).Value = "Hortaliza"
